$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 121, shifting rows 121:153 down to 122:154.
$ws.Rows(121).Insert()

# Populate the newly inserted row 121 with a new data record.
$ws.Cells.Item(121, 1).Value = 8
$ws.Cells.Item(121, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(121, 3).Value = "Coquimbo"
$ws.Cells.Item(121, 4).Value = 44463
$ws.Cells.Item(121, 5).Value = 4
$ws.Cells.Item(121, 6).Value = 100112012
$ws.Cells.Item(121, 7).Value = "Espinaca"
$ws.Cells.Item(121, 8).Value = "Sin especificar"
$ws.Cells.Item(121, 9).Value = "Primera"
$ws.Cells.Item(121, 10).Value = 3400
$ws.Cells.Item(121, 11).Value = 400
$ws.Cells.Item(121, 12).Value = 500
$ws.Cells.Item(121, 13).Value = 450
$ws.Cells.Item(121, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(121, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(121, 16).Value = 900
$ws.Cells.Item(121, 17).Value = 0.5
$ws.Cells.Item(121, 18).Value = "Hortaliza"
